# SetoresRelevantes.xlsx — add column B "SECRETARIA/COORDENADORIA RESP"
# mapping each sector (column A) to its parent Secretaria/Coordenadoria,
# fixing duplicated sector names across departments.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Header for the new column B
$ws.Range("B1").Value = "SECRETARIA/COORDENADORIA RESP"

# Parent department labels
$secGS      = "Secretaria de Gestão de Serviços"
$coordInfra = "Coordenadoria de InfraEstrutura Predial"
$coordSeg   = "Coordenadoria de Segurança , Transporte e Apoio Administrativo"

# Rows 2-3: SECGS / GABGS -> Secretaria de Gestão de Serviços
$ws.Range("B2").Value = $secGS
$ws.Range("B3").Value = $secGS

# Rows 4-13: CIP, SAPRE, CAA, SMOEP, SMOP, SMOI, SGACI, SMIC, SMIN, SOP
#   -> Coordenadoria de InfraEstrutura Predial
$ws.Range("B4").Value = $coordInfra
$ws.Range("B5").Value = $coordInfra
$ws.Range("B6").Value = $coordInfra
$ws.Range("B7").Value = $coordInfra
$ws.Range("B8").Value = $coordInfra
$ws.Range("B9").Value = $coordInfra
$ws.Range("B10").Value = $coordInfra
$ws.Range("B11").Value = $coordInfra
$ws.Range("B12").Value = $coordInfra
$ws.Range("B13").Value = $coordInfra

# Rows 14-23: CSTA, SEXP, ST, ASSISEG, SESEG, SST, SMI, ASG, SAFI, COGSA
#   -> Coordenadoria de Segurança , Transporte e Apoio Administrativo
$ws.Range("B14").Value = $coordSeg
$ws.Range("B15").Value = $coordSeg
$ws.Range("B16").Value = $coordSeg
$ws.Range("B17").Value = $coordSeg
$ws.Range("B18").Value = $coordSeg
$ws.Range("B19").Value = $coordSeg
$ws.Range("B20").Value = $coordSeg
$ws.Range("B21").Value = $coordSeg
$ws.Range("B22").Value = $coordSeg
$ws.Range("B23").Value = $coordSeg

# Column widths (stored OOXML width = ColumnWidth input + 0.8333; compensate
# so the saved file carries width=26.5 / 50.875 as in the target workbook)
$ws.Columns.Item(1).ColumnWidth = 25.666666666666668
$ws.Columns.Item(2).ColumnWidth = 50.041666666666664

# Match the author's final selection state
$ws.Range("E1:E3").Select()
